# ------------------------------------------------------------------
# Add a new "Player Info" sheet (scraped player attributes) ahead of
# the existing "ODI Batting" sheet, and tweak ODI Batting's
# MATCH_CARD_LINK column into a plain MATCH_CODE column.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# The workbook currently holds a single sheet: "ODI Batting".
$odiBatting = $wb.Worksheets.Item(1)

# --- Step 1: capture the existing ODI Batting data before we repurpose it ---
# (NOTE: ".Value" getter is unreliable in this host -- use ".Value2" to read.)
$headers = @("MATCH_NUMBER","INNING_NUMBER","MATCH_DATE","MATCH_CARD_LINK","MATCH_INNING","OPPONENT","VENUE","DISMISSAL","RUNS_SCORED","BALLS_FACED")
$rowValues = @()
for ($col = 1; $col -le 10; $col++) {
    $rowValues += ,$odiBatting.Cells.Item(2, $col).Value2
}

# --- Step 2: insert a brand-new sheet right after the original sheet.  ---
# --- The *existing* sheet object becomes "Player Info" (keeping        ---
# --- sheetId 1 / rId1), while the freshly added sheet becomes          ---
# --- "ODI Batting" with sheetId 2 / rId2 and carries the original data ---
$newOdiBatting = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $odiBatting)

# Free up the "ODI Batting" name on the original sheet before claiming it
# on the new sheet (Excel disallows duplicate sheet names).
$odiBatting.Name = "Player Info"
$newOdiBatting.Name = "ODI Batting"

for ($col = 1; $col -le 10; $col++) {
    $newOdiBatting.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$newOdiBatting.Range("D1").Value = "MATCH_CODE"

# Re-apply the bold / bordered / centered header look used throughout the
# workbook (matches cellXfs index 1 in the original styles.xml).
$headerRange = $newOdiBatting.Range("A1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous (thin box border)

for ($col = 1; $col -le 10; $col++) {
    $val = $rowValues[$col - 1]
    $cell = $newOdiBatting.Cells.Item(2, $col)
    if ($col -eq 4) {
        $val = "4447"
    }
    if ($val -match '^[0-9]+$') {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

# --- Step 3: rewrite the original sheet object's content into "Player Info" ---
$odiBatting.Range("E1:J2").Clear()

$odiBatting.Range("A1").Value = "ID"
$odiBatting.Range("B1").Value = "NAME"
$odiBatting.Range("C1").Value = "BATTING_HAND"
$odiBatting.Range("D1").Value = "BOWL_STYLE"

$odiBatting.Range("A2").NumberFormat = "@"
$odiBatting.Range("A2").Value = "5402"
$odiBatting.Range("A2").Style = "Normal"
$odiBatting.Range("B2").Value = "Jahmar Neville Hamilton"
$odiBatting.Range("C2").Value = "Right Handed"
$odiBatting.Range("D2").Value = "Does Not Bowl | Unknown"

Write-Host "Workbook restructured: Player Info + ODI Batting"
